$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header rows: each data block begins with a header row whose three
# cells (A/B/C) identify the source edf file; they now carry distinct
# column headers instead of the generic "A"/"B" placeholder.
$headerRows = @(1, 14, 27, 40)
foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 1).Value = "ABC"
    $ws.Cells.Item($r, 2).Value = "DEF"
    $ws.Cells.Item($r, 3).Value = "GHI"
}

# Normalize the workbook's font to Calibri (was MS PGothic) by updating the
# workbook's Normal style, so every cell (current and future) picks it up
# without stamping explicit per-cell styles on blank/unused rows.
$wb.Styles.Item("Normal").Font.Name = "Calibri"
